# Apply the edit described by the diff:
# Insert a new row at row 728 with data (2026/01/26, 月, 16, 163),
# shifting all subsequent rows (728-769) down to (729-770).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 728, shifting existing rows down.
$ws.Rows("728:728").Insert()

# Populate the newly inserted row 728 with the new data.
# Column A holds a date-looking string ("2026/01/26") that must stay a
# plain text value (matching the other date cells in this sheet, which
# are plain strings rather than date serials). Temporarily force a Text
# number format so the value isn't auto-parsed as a date, then reset the
# cell style back to Normal so no stray number-format style is left
# behind (matching the unstyled cells elsewhere in the column).
$ws.Cells.Item(728, 1).NumberFormat = "@"
$ws.Cells.Item(728, 1).Value = "2026/01/26"
$ws.Cells.Item(728, 1).Style = "Normal"
$ws.Cells.Item(728, 2).Value = "月"
$ws.Cells.Item(728, 3).Value = 16
$ws.Cells.Item(728, 4).Value = 163
